# product_in_stock.xlsx — update stock quantity and cursor position
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The quantity of the first product (row 2, column B) changed from 10 to 0
$ws.Range("B2").Value = 0

# The active cell / selection moved from C9 to B2
$ws.Activate()
$ws.Range("B2").Select()
